$d = $word.ActiveDocument

# Locate the "Requisitos" entry paragraph ("LOB1053: Física III (Requisito
# fraco)") with Find. The three paragraphs that immediately follow it must
# be removed in their entirety:
#   1) the blank paragraph
#   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3) "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github
#       pages. Original theme under Creative Commons Attribution"
# leaving the trailing blank paragraph + page-break paragraph untouched.

$findRange = $d.Content
$findRange.Find.Execute("LOB1053: Física III (Requisito fraco)", $true, $false,
                         $false, $false, $false, $true, 1, $false, "", 0)

if ($findRange.Find.Found) {
    # Map the found range back to its paragraph index so we can reliably
    # walk to the following paragraphs by index (Paragraph.Next proved
    # unreliable across this anchor).
    $anchorIndex = -1
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Start -le $findRange.Start -and $p.Range.End -ge $findRange.End) {
            $anchorIndex = $i
        }
    }

    if ($anchorIndex -gt 0) {
        $firstToRemove = $d.Paragraphs.Item($anchorIndex + 1)
        $lastToRemove = $d.Paragraphs.Item($anchorIndex + 3)

        $deleteRange = $d.Range($firstToRemove.Range.Start, $lastToRemove.Range.End)
        $deleteRange.Delete()
    }
}

Write-Output "done"
